$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H3").Value = 2.66
$ws.Range("I3").Value = 3.65
$ws.Range("J3").Value = 2.82
$ws.Range("K3").Value = 950
$ws.Range("G4").Value = 1.44
$ws.Range("F5").Value = 5.2
$ws.Range("G5").Value = 5.7
$ws.Range("H5").Value = 1.83
$ws.Range("I5").Value = 1.85
$ws.Range("J5").Value = 3.65
$ws.Range("K5").Value = 3.85
$ws.Range("G6").Value = 2.2
$ws.Range("J6").Value = 3.2
$ws.Range("Q6").Value = 1.46
$ws.Range("I7").Value = 2
$ws.Range("P7").Value = 1.81
$ws.Range("Q7").Value = 1.74
$ws.Range("Q10").Value = 1.84
$ws.Range("F11").Value = 4.2
$ws.Range("G11").Value = 5.4
$ws.Range("H11").Value = 1.69
$ws.Range("I11").Value = 1.93
$ws.Range("J11").Value = 4.1
$ws.Range("K11").Value = 5.7
$ws.Range("P11").Value = 2.66
$ws.Range("Q11").Value = 1.42
$ws.Range("F12").Value = 1.7
$ws.Range("P12").Value = 2.06
$ws.Range("Q12").Value = 1.59
$ws.Range("F13").Value = 3.55
$ws.Range("G13").Value = 4.2
$ws.Range("H13").Value = 1.86
$ws.Range("I13").Value = 2
$ws.Range("P13").Value = 2.48
$ws.Range("Q13").Value = 1.37
$ws.Range("I14").Value = 3.2
$ws.Range("F15").Value = 1.97
$ws.Range("H15").Value = 1.67
$ws.Range("I15").Value = 2.1
$ws.Range("J15").Value = 1.91
$ws.Range("P15").Value = 3.45
$ws.Range("Q15").Value = 1.28
$ws.Range("G19").Value = 4.9
$ws.Range("H19").Value = 1.97
$ws.Range("K19").Value = 3.5
$ws.Range("I20").Value = 5.5
$ws.Range("K20").Value = 3.8
$ws.Range("R20").Value = 1.25
$ws.Range("U20").Value = 1.79
$ws.Range("Z20").Value = 44
$ws.Range("AD20").Value = 23
$ws.Range("AE20").Value = 100
$ws.Range("F21").Value = 2.56
$ws.Range("G21").Value = 2.62
$ws.Range("H21").Value = 3.2
$ws.Range("I21").Value = 3.3
$ws.Range("O21").Value = 1.44
$ws.Range("Y21").Value = 10.5
$ws.Range("Z21").Value = 21
$ws.Range("AA21").Value = 60
$ws.Range("AE21").Value = 44
$ws.Range("AI21").Value = 60
$ws.Range("AJ21").Value = 40
$ws.Range("AO21").Value = 50
$ws.Range("G26").Value = 3.5
$ws.Range("F29").Value = 3.25
$ws.Range("G29").Value = 3.85
$ws.Range("H29").Value = 2.36
$ws.Range("I29").Value = 2.86
$ws.Range("J29").Value = 2.92
$ws.Range("K29").Value = 3.75
$ws.Range("P29").Value = 1.53
$ws.Range("Q29").Value = 2.5
